# Append the new record row (row 87) to the "AYKO" sheet, matching the
# automated map-export format used for this workbook: every column is
# stored as literal text except "Attachments" (col I), which is numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 87

# (column index, text value) pairs for every text-typed column in the row.
$textCells = @(
    @(1,  "-521"),
    @(2,  "7/21/2025"),
    @(3,  "Gregario Laferrere 6572"),
    @(4,  "9"),
    @(5,  "808430950"),
    @(6,  "AYKO"),
    @(7,  "Pendiente"),
    @(8,  "Entre edificio 12 y edificio 13 columna corroida en base tambalea"),
    @(10, "Cambio"),
    @(11, "Sin equipos"),
    @(12, "Pasante"),
    @(13, ""),
    @(14, ""),
    @(15, "No ubicado"),
    @(16, "No clasificado, consultar con mantenimiento")
)

foreach ($pair in $textCells) {
    $col = $pair[0]
    $val = $pair[1]
    $cell = $ws.Cells.Item($row, $col)
    # Force text interpretation so numeric-looking strings (e.g. "-521",
    # "9", "808430950") and the literal date string stay as text instead
    # of being coerced to numbers/dates, then drop back to the default
    # "Normal" style so no explicit style id is left on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# "Attachments" column is numeric.
$ws.Cells.Item($row, 9).Value = 1
